$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Majorelle Magdy"
$ws.Cells.Item(3, 7).Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Cells.Item(4, 7).Value = "Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(5, 7).Value = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(6, 7).Value = "Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(7, 7).Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(8, 7).Value = "Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Administrator, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Asmaa Reda"
$ws.Cells.Item(9, 7).Value = "Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(10, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Sara Wael, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel"
$ws.Cells.Item(11, 7).Value = "Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda"
$ws.Cells.Item(13, 7).Value = "Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad, D Wessam Atef"
$ws.Cells.Item(15, 7).Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(16, 7).Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Cells.Item(17, 7).Value = "Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Basma Hamed"
$ws.Cells.Item(19, 7).Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Cells.Item(22, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(24, 7).Value = "Dr. Aya Emad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Atef, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Remon, Dr. Youstina Magdy"
$ws.Cells.Item(25, 7).Value = "Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy"
$ws.Cells.Item(27, 7).Value = "Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Remon"
$ws.Cells.Item(28, 7).Value = "Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon"
$ws.Cells.Item(29, 7).Value = "Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Remon"
$ws.Cells.Item(30, 7).Value = "Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(31, 7).Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Cells.Item(32, 7).Value = "Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(33, 7).Value = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(34, 7).Value = "Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(35, 7).Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(36, 7).Value = "Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Administrator, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Asmaa Reda"
$ws.Cells.Item(37, 7).Value = "Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Gehan Adel"
$ws.Cells.Item(38, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Sara Wael, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel"
$ws.Cells.Item(39, 7).Value = "Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda"
$ws.Cells.Item(41, 7).Value = "Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad, D Wessam Atef"
$ws.Cells.Item(43, 7).Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Cells.Item(44, 7).Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Cells.Item(45, 7).Value = "Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Basma Hamed"
$ws.Cells.Item(47, 7).Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Cells.Item(50, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(52, 7).Value = "Dr. Aya Emad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Atef, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Remon, Dr. Youstina Magdy"
$ws.Cells.Item(53, 7).Value = "Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy"
$ws.Cells.Item(55, 7).Value = "Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Remon"
$ws.Cells.Item(56, 7).Value = "Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon"
$ws.Cells.Item(57, 7).Value = "Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Remon"
